$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) In-place item-list edits (no row-count change yet).
# ------------------------------------------------------------------
# weights - iron total: 804 -> 824
$ws.Range("B13").Value = 824

# row21: "wrist straps"/10 -> "dip belt"/40
$ws.Range("A21").Value = "dip belt"
$ws.Range("B21").Value = 40

# row22: "chalk"/8 -> "chalk (x2)"/16
$ws.Range("A22").Value = "chalk (x2)"
$ws.Range("B22").Value = 16

# ------------------------------------------------------------------
# 2) Insert a new row at 16 for "dumbbells - 90lb" / ".180."
#    (this shifts the row16.. block down by one row)
# ------------------------------------------------------------------
$ws.Rows(16).Insert()
$ws.Range("A17").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = "dumbbells - 90lb"
$ws.Range("B16").Value = ".180."
$ws.Range("B16").HorizontalAlignment = -4152  # xlRight

# ------------------------------------------------------------------
# 3) Insert two more rows right after the "dip belt" row (now row 22)
#    to hold "lifting straps" and "wrist wraps" ahead of "chalk (x2)"
#    (now pushed down to row 25).
# ------------------------------------------------------------------
$ws.Rows(23).Insert()
$ws.Rows(23).Insert()

$ws.Range("A22").Copy() | Out-Null
$ws.Range("A24").PasteSpecial(-4122) | Out-Null
$ws.Range("B22").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Range("A24").Value = "lifting straps"
$ws.Range("B24").Value = 8

$ws.Range("A22").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null
$ws.Range("B22").Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$ws.Range("A23").Value = "wrist wraps"
$ws.Range("B23").Value = 10

# ------------------------------------------------------------------
# 4) Fix up the Total row: formula now spans to B29, value 3132.
#    (after the 3 row-inserts above, the old row 28 "Total" is now row 31)
# ------------------------------------------------------------------
$ws.Range("B31").Formula = "=SUM(B2:B29)"

# ------------------------------------------------------------------
# 6) Add 3 extra blank rows at the end of the sheet (A1:E43 -> A1:E46)
# ------------------------------------------------------------------
$ws.Range("A43:E43").Copy() | Out-Null
$ws.Range("A44:E46").PasteSpecial(-4122) | Out-Null

# ------------------------------------------------------------------
# 7) Update the comment on B13 (the weights-iron breakdown).
# ------------------------------------------------------------------
$comment = $ws.Range("B13").Comment
$comment.Text("$352 - 45 x 4" + [char]10 + "$322 - 25 x 8" + [char]10 + "$  86 - 10 x 4" + [char]10 + "$  46 - 5   x 4" + [char]10 + "$  28 - 2.5 x 4")

# ------------------------------------------------------------------
# 8) Selection + window size cosmetic changes.
# ------------------------------------------------------------------
$ws.Range("C30").Select()
$excel.ActiveWindow.Width = 23700
$excel.ActiveWindow.Height = 14085
